$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark at the very start of the document.
# ---------------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
    # no-op if it somehow doesn't exist
}

# ---------------------------------------------------------------------------
# 2) Rework the "meer dan vier keer te laat komt" sentence:
#    - insert "zonder geldige reden " before "te laat komt, zal"
#    - drop ", wanneer diegene geen geldige reden heeft om te laat te komen."
#      and replace with "."
#    - insert " van afwezigheid" after "(Wanneer een reden"
# ---------------------------------------------------------------------------
$rInsert1 = $d.Content
$rInsert1.Find.Execute("keer te laat komt, zal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $rInsert1.Start + 5
$keerRange = $d.Range($rInsert1.Start, $splitPoint)
$keerRange.InsertAfter("zonder geldige reden ")

$rReplace = $d.Content
$rReplace.Find.Execute(", wanneer diegene geen geldige reden heeft om te laat te komen.", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2)

$rInsert2 = $d.Content
$rInsert2.Find.Execute("(Wanneer een reden", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rInsert2.InsertAfter(" van afwezigheid")

# Re-split the paragraph's single merged run back into the distinct runs seen
# in the target document (text content is already correct at this point; we
# only need to re-introduce the original run boundaries). Toggling Bold on
# and back off on a sub-range forces Word to break the run at that boundary
# without altering any visible formatting.
$rParagraph = $d.Content
$rParagraph.Find.Execute("Wanneer iemand meer dan vier keer zonder geldige reden te laat komt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pStart = $rParagraph.Start

$pRange = $d.Range($pStart, $pStart)
$pRange.Expand(4) | Out-Null
$pEnd = $pRange.End

$splitOffsets = @(6, 29, 34, 55, 72, 104, 105, 124, 140)
foreach ($off in $splitOffsets) {
    $pos = $pStart + $off
    $splitRange = $d.Range($pos, $pEnd)
    $splitRange.Bold = 1
    $splitRange.Bold = 0
}

# ---------------------------------------------------------------------------
# 3) Fix "zsm" -> "z.s.m." (dropping the spell-check proofErr wrapper) and
#    plant a fresh "_GoBack" bookmark right after the corrected word.
# ---------------------------------------------------------------------------
$rZsm = $d.Content
$rZsm.Find.Execute("Marco probeert elke dag om 09.30 aanwezig te zijn, mocht dit niet lukken dan laat hij dit de rest van de groep weten en probeert hij zsm op school aanwezig te zijn.", $true, $false, $false, $false, $false, $true, 1, $false, "Marco probeert elke dag om 09.30 aanwezig te zijn, mocht dit niet lukken dan laat hij dit de rest van de groep weten en probeert hij z.s.m. op school aanwezig te zijn.", 2)

$rZsmParagraph = $d.Content
$rZsmParagraph.Find.Execute("Marco probeert elke dag om 09.30 aanwezig te zijn, mocht dit niet lukken dan laat hij dit de rest van de groep weten en probeert hij z.s.m. op school aanwezig te zijn.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$zsmPStart = $rZsmParagraph.Start
$zsmPEnd = $rZsmParagraph.End

$runA_len = 133   # length of "Marco probeert ... probeert hij "
$runB_len = 139   # offset where "z.s.m." ends / trailing text begins

$posRunB = $zsmPStart + $runA_len
$posRunC = $zsmPStart + $runB_len

$splitB = $d.Range($posRunB, $zsmPEnd)
$splitB.Bold = 1
$splitB.Bold = 0

$splitC = $d.Range($posRunC, $zsmPEnd)
$splitC.Bold = 1
$splitC.Bold = 0

$newBookmarkRange = $d.Range($posRunC, $posRunC)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
